$d = $word.ActiveDocument

# Helper: replace the full text of a paragraph (drops the trailing paragraph
# mark from the range so we don't delete the paragraph itself) and collapses
# the paragraph down to a single run, which also removes any <w:br/> that
# used to sit inside it.
function Set-ParaText($para, [string]$newText) {
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

# --- Paragraph 1: title line - update date and paper title (two runs,
#     separated by a manual line break that must be preserved) ---
$d.Content.Find.Execute("25.12.24", $true, $false, $false, $false, $false, $true, 1, $false, "23.12.24", 2) | Out-Null
$d.Content.Find.Execute("Vision language models are blind", $true, $false, $false, $false, $false, $true, 1, $false, "T-FREE: Tokenizer-Free Generative LLMs via Sparse Representations for Memory-Efficient Embeddings", 2) | Out-Null

# --- Paragraph 2: intro paragraph rewritten entirely ---
Set-ParaText $d.Paragraphs.Item(2) "שוב חוזרים לנושא הטוקנייזרים - מתברר שהוא יותר חם ממה שחשבתי. נתקלתי במאמר המעניין שיטה נוספת לטוקניזציה המבוססת על פונקציה האש n-grams. השיטה המוצעת באה להתמודד עם גודל העצום של המילון מלווה כל מודל שפה גדול (עשרות אלפי טוקנים לכל הפחות) וגם טוקנים דומים מאוד מבחינת האותיות  האותיות שמצריכות אמבדינגים שונים שזה לא יעיל (לטענת המחברים). "

# --- Paragraph 3: rewritten, and its trailing <w:br/> is removed ---
Set-ParaText $d.Paragraphs.Item(3) "המחברים מנסים שיטת טוקניזציה שה-encoding שלה המורכב משלבים הבאים:"

# --- Paragraph 4: rewritten, and its trailing <w:br/> is removed ---
Set-ParaText $d.Paragraphs.Item(4) "פירוק של טקסט למה שהם קוראים טוקנים כאשר ב-T-FREE טוקנים אלו הם בעצם מילים"

# --- Paragraph 5: rewritten, and its trailing <w:br/> is removed ---
Set-ParaText $d.Paragraphs.Item(5) "כל מילה מחולקת לסדרה של 3-grams לא זרים למשל מילה hello מיוצגת על ידי חמישה 3-grams הבאים: {_He, Hel, ell, llo, lo_}. מספר 3-grams בייצוג הזה בדרך כלל מספר n-grams במילה שווה למספר האותיות במילה"

# --- Paragraph 6: rewritten, and its trailing <w:br/> is removed ---
Set-ParaText $d.Paragraphs.Item(6) "מקודדים כל 3-gram עם m פונקציות האש שכל אחת מהם מקבלת v ערכים אפשריים כאשר v הינו אחד הייפר-הפרמטרים של השיטה."

# --- Paragraph 7: rewritten, and its trailing <w:br/> is removed ---
Set-ParaText $d.Paragraphs.Item(7) "כך כל מילה מקודדת על ידי n*m מספרים בין 0 ל-v כאשר n הינו אורך המילה (מספר אותיות). ייצוג המילה הוא ממוצע (ועיגול) של כל nm ערכים האלו."

# --- Paragraph 8: rewritten (no break involved) ---
Set-ParaText $d.Paragraphs.Item(8) "כל ערך בין 0 ל-v מקודד על ידי וקטור נלמד כאשר v וקטורים אלו למעשה מהווים את המילון של השיטה"

# --- Paragraph 9: rewritten (no break involved) ---
Set-ParaText $d.Paragraphs.Item(9) "שלב האימון והפענוח (כלומר גנרוט של מילים) נראים קצת יותר מורכבים. קודם כל באימון המטרה היא לחזות את nm האשים של 3-grams של המילה הבאה. כלומר במקום בעיית multi-class בפענוח של הטוקניזציה הרגילה (חיזוי של טוקן ממילון הטוקנים) יש לנו כאן בעיית multi-label כאשר אנו חוזים n*m האשים. שימו לב ש n תלוי באורך המילה כלומר יש לנו מספר ״לייבלים״ שונה לפי אורך המילה."

# --- Paragraph 10 (previously the arXiv link paragraph): replaced with the
#     decoding-discussion paragraph, and two brand-new paragraphs (the
#     closing remark and the new arXiv link) are appended after it ---
Set-ParaText $d.Paragraphs.Item(10) "הפענוח לא ממש ברור לי האמת. כאשר אנו רוצים לחזות את המילה הבאה אנו קודם כל מחשבים את כל ההאשים עבור כל המילים האפשריות (זה די הרבה כי לכל מילה יש גם את כל ההטיות שלה לכל הפחות ובנוסף מילים בעלות אורכים שונים מקודדים עם מספר n*m שונה של האשים). לאחר מכן בוחרים את המילה המיוצגת על יד האשים בעלי ״ההסתברות הגבוהה ביותר״. נזכור שהמודל חוזה הסתברות של כל ערך של האש מ 1 עד v (גודל המילון) ולא לגמרי ברור איך נבחרת קבוצת האשים בעלת הסתברות הגבוהה ביותר."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
Set-ParaText $d.Paragraphs.Item($d.Paragraphs.Count) "בקיצור מאמר נחמד אבל לא ברור לי העניין עם הפענוח…"

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
Set-ParaText $d.Paragraphs.Item($d.Paragraphs.Count) "https://arxiv.org/abs/2406.19223"

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
